$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Originally row 2 held the single FAPs/Rspo2/Rnf43 -> MuSCs summary line.
# The new TPM numbers split this into two target-cluster rows: a new "ECs"
# row and the pre-existing "MuSCs" row (now with recomputed specificity
# figures because a second target cluster exists). Read the current row 2
# contents first (A:L are identical between both rows) so they can be
# carried over to the new row 3 without disturbing row 2's own cells.

$commonCols = 1..12
$orig = @{}
foreach ($col in $commonCols) {
    $orig[$col] = $ws.Cells.Item(2, $col).Value2
}

# Row 3 = the original MuSCs row, pushed down, with updated M:T values.
foreach ($col in $commonCols) {
    $ws.Cells.Item(3, $col).Value = $orig[$col]
}
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 13).Value = 0.04441766666666667
$ws.Cells.Item(3, 14).Value = 0.133253
$ws.Cells.Item(3, 15).Value = 0.8637760261363342
$ws.Cells.Item(3, 16).Value = 0.8637760261363342
$ws.Cells.Item(3, 17).Value = 0.007579016075111112
$ws.Cells.Item(3, 18).Value = 0.06821114467600001
$ws.Cells.Item(3, 19).Value = 0.8637760261363342
$ws.Cells.Item(3, 20).Value = 0.8637760261363342

# Row 2 = new ECs target-cluster row: only the target-cluster label (D2) and
# the M:T edge-weight figures actually change; A2:L2 stay untouched.
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 13).Value = 0.007005
$ws.Cells.Item(2, 14).Value = 0.021015
$ws.Cells.Item(2, 15).Value = 0.1362239738636658
$ws.Cells.Item(2, 16).Value = 0.1362239738636658
$ws.Cells.Item(2, 17).Value = 0.00119526782
$ws.Cells.Item(2, 18).Value = 0.01075741038
$ws.Cells.Item(2, 19).Value = 0.1362239738636658
$ws.Cells.Item(2, 20).Value = 0.1362239738636658
